$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "'581-08055C103J"
$ws.Range("B5").Value = "'MOUSER"
$ws.Range("C5").Value = "'CAPACITOR, 0.01uF, 50V"

$ws.Range("A5:C5").Interior.Color = 65535

$ws.Range("C10").Select()
